$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(85, 1).Value = 45515
$ws.Cells.Item(85, 2).Value = 746.4783276449
$ws.Cells.Item(85, 3).Value = 177.180722063
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 9).Value = 253.473029667
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0.06583642968800001
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0
$ws.Cells.Item(85, 14).Value = 96.38300855904001
$ws.Cells.Item(85, 15).Value = 51.145575918
$ws.Cells.Item(85, 16).Value = 0
$ws.Cells.Item(85, 17).Value = 0.0000018552
$ws.Cells.Item(85, 18).Value = 0
$ws.Cells.Item(85, 19).Value = 0
$ws.Cells.Item(85, 20).Value = 0
$ws.Cells.Item(85, 21).Value = 257.6490425180494
$ws.Cells.Item(85, 23).Value = 0
$ws.Cells.Item(85, 24).Value = 0
$ws.Cells.Item(85, 25).Value = 0
$ws.Cells.Item(85, 26).Value = 205.87864775833

$ws.Cells.Item(86, 1).Value = 45516
$ws.Cells.Item(86, 2).Value = 754.5397090904
$ws.Cells.Item(86, 3).Value = 188.754345605
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 9).Value = 261.923921099
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0
$ws.Cells.Item(86, 14).Value = 104.82531587808
$ws.Cells.Item(86, 15).Value = 52.71457384200001
$ws.Cells.Item(86, 16).Value = 0
$ws.Cells.Item(86, 17).Value = 0.0000019728
$ws.Cells.Item(86, 18).Value = 0
$ws.Cells.Item(86, 19).Value = 0
$ws.Cells.Item(86, 20).Value = 0
$ws.Cells.Item(86, 23).Value = 0
$ws.Cells.Item(86, 24).Value = 0
$ws.Cells.Item(86, 25).Value = 0
$ws.Cells.Item(86, 26).Value = 194.95532384998

$ws.Cells.Item(87, 1).Value = 45517
$ws.Cells.Item(87, 2).Value = 770.3116896865
$ws.Cells.Item(87, 3).Value = 187.377325694
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 9).Value = 262.568480615
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0.001230713055
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0
$ws.Cells.Item(87, 14).Value = 106.66788295168
$ws.Cells.Item(87, 15).Value = 53.32493170200001
$ws.Cells.Item(87, 16).Value = 0
$ws.Cells.Item(87, 17).Value = 0.0000019632
$ws.Cells.Item(87, 18).Value = 0
$ws.Cells.Item(87, 19).Value = 0
$ws.Cells.Item(87, 20).Value = 0
$ws.Cells.Item(87, 21).Value = 266.8599318235606
$ws.Cells.Item(87, 23).Value = 0
$ws.Cells.Item(87, 24).Value = 0
$ws.Cells.Item(87, 25).Value = 0
$ws.Cells.Item(87, 26).Value = 202.820117063992

$ws.Cells.Item(88, 1).Value = 45518
$ws.Cells.Item(88, 2).Value = 746.1070756329
$ws.Cells.Item(88, 3).Value = 184.5352287075
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 9).Value = 257.304577901
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0.00249393912
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0
$ws.Cells.Item(88, 14).Value = 103.45176587776
$ws.Cells.Item(88, 15).Value = 53.243550654
$ws.Cells.Item(88, 16).Value = 0
$ws.Cells.Item(88, 17).Value = 0.00000192
$ws.Cells.Item(88, 18).Value = 0
$ws.Cells.Item(88, 19).Value = 0
$ws.Cells.Item(88, 20).Value = 0
$ws.Cells.Item(88, 21).Value = 256.7535393911247
$ws.Cells.Item(88, 23).Value = 0
$ws.Cells.Item(88, 24).Value = 0
$ws.Cells.Item(88, 25).Value = 0
$ws.Cells.Item(88, 26).Value = 202.924148720262

$ws.Cells.Item(89, 1).Value = 45519
$ws.Cells.Item(89, 2).Value = 731.5833663566
$ws.Cells.Item(89, 3).Value = 178.1867925015
$ws.Cells.Item(89, 4).Value = 0
$ws.Cells.Item(89, 5).Value = 0
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 9).Value = 255.281377198
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0.00373517109
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = 0
$ws.Cells.Item(89, 14).Value = 92.89888172895999
$ws.Cells.Item(89, 15).Value = 52.816300152
$ws.Cells.Item(89, 16).Value = 0
$ws.Cells.Item(89, 17).Value = 0.0000017856
$ws.Cells.Item(89, 18).Value = 0
$ws.Cells.Item(89, 19).Value = 0
$ws.Cells.Item(89, 20).Value = 0
$ws.Cells.Item(89, 21).Value = 250.9967335751802
$ws.Cells.Item(89, 23).Value = 0
$ws.Cells.Item(89, 24).Value = 0
$ws.Cells.Item(89, 25).Value = 0
$ws.Cells.Item(89, 26).Value = 195.496288462584

$ws.Cells.Item(90, 1).Value = 45520
$ws.Cells.Item(90, 2).Value = 748.5381406060001
$ws.Cells.Item(90, 3).Value = 179.7704347355
$ws.Cells.Item(90, 4).Value = 0
$ws.Cells.Item(90, 5).Value = 0
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 9).Value = 248.943208624
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 0.00512940546
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 13).Value = 0
$ws.Cells.Item(90, 14).Value = 87.94070124000001
$ws.Cells.Item(90, 15).Value = 52.86716330700001
$ws.Cells.Item(90, 16).Value = 0
$ws.Cells.Item(90, 17).Value = 0.0000017232
$ws.Cells.Item(90, 18).Value = 0
$ws.Cells.Item(90, 19).Value = 0
$ws.Cells.Item(90, 20).Value = 0
$ws.Cells.Item(90, 21).Value = 248.4381532125382
$ws.Cells.Item(90, 23).Value = 0
$ws.Cells.Item(90, 24).Value = 0
$ws.Cells.Item(90, 25).Value = 0
$ws.Cells.Item(90, 26).Value = 191.792761499372

$ws.Cells.Item(91, 1).Value = 45521
$ws.Cells.Item(91, 2).Value = 756.3877049789
$ws.Cells.Item(91, 3).Value = 181.2805804385
$ws.Cells.Item(91, 4).Value = 0
$ws.Cells.Item(91, 5).Value = 0
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 9).Value = 254.117589183
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 0.006450007425
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).Value = 0
$ws.Cells.Item(91, 14).Value = 106.63438173216
$ws.Cells.Item(91, 15).Value = 55.094969496
$ws.Cells.Item(91, 16).Value = 0
$ws.Cells.Item(91, 17).Value = 0.0000017976
$ws.Cells.Item(91, 18).Value = 0
$ws.Cells.Item(91, 19).Value = 0
$ws.Cells.Item(91, 20).Value = 0
$ws.Cells.Item(91, 21).Value = 268.1392220048816
$ws.Cells.Item(91, 23).Value = 0
$ws.Cells.Item(91, 24).Value = 0
$ws.Cells.Item(91, 25).Value = 0
$ws.Cells.Item(91, 26).Value = 192.416951436992

$ws.Cells.Item(92, 1).Value = 45522
$ws.Cells.Item(92, 3).Value = 181.1169466525
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0.00775148409
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(92, 14).Value = 102.31272441408
$ws.Cells.Item(92, 16).Value = 0
$ws.Cells.Item(92, 17).Value = 0.0000017496
$ws.Cells.Item(92, 20).Value = 0
$ws.Cells.Item(92, 24).Value = 0
$ws.Cells.Item(92, 26).Value = 179.621057715782

$ws.Range("A84").Copy()
$ws.Range("A85:A92").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0
